$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "Cable type list" sheet between "Connector List" and "Cable Order" ---
$afterSheet = $wb.Worksheets.Item("Connector List")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "Cable type list"

# Header row
$newSheet.Range("A1").Value = "Cable Signal Type"
$newSheet.Range("B1").Value = "Length to be purchased (mm)"

# Data rows
$newSheet.Range("A2").Value = "CAN"
$newSheet.Range("B2").Value = 10600

$newSheet.Range("A3").Value = "Electric (Encoder)"
$newSheet.Range("B3").Value = 4200

$newSheet.Range("A4").Value = "Electric (Hall)"
$newSheet.Range("B4").Value = 4200

$newSheet.Range("A5").Value = "Electric (Motor)"
$newSheet.Range("B5").Value = 4200

$newSheet.Range("A6").Value = "Ethernet"
$newSheet.Range("B6").Value = 10890

$newSheet.Range("A7").Value = "USB 2.0"
$newSheet.Range("B7").Value = 1050

# --- 2. Rename "10pinF / SH-Cut" references on the "Cable List" sheet ---
$cableList = $wb.Worksheets.Item("Cable List")
$cableList.Range("F13").Value = "10pinF(CAN) / SH-Cut"
$cableList.Range("F16").Value = "10pinF(CAN) / SH-Cut"

# --- 3. Update the selection on "Cable Order 2" (was F2, now D6) ---
$cableOrder2 = $wb.Worksheets.Item("Cable Order 2")
$cableOrder2.Activate()
$cableOrder2.Range("D6").Select()

# --- 4. Restore active sheet to "Cable List" and update its selection (was H7, now F10) ---
$cableList.Activate()
$cableList.Range("F10").Select()
